# Apply the authors edit: rows 4-6 are re-sorted (rotated) and a new row 7 is appended
# (matches the canonical OOXML diff: new-row4 = old-row5, new-row5 = old-row6, new-row6 = old-row4, plus a new row 7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("I4").NumberFormat = "@"
$ws.Range("A4").Value = 111803593
$ws.Range("B4").Value = 85188
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 1980
$ws.Range("F4").Value = "Porslinsblå spindling"
$ws.Range("G4").Value = "Cortinarius cumatilis"
$ws.Range("H4").Value = "Fr."
$ws.Range("I4").Value = "20"
$ws.Range("Q4").Value = 695812.1972037496
$ws.Range("R4").Value = 6553542.22335465
$ws.Range("S4").Value = 2
$ws.Range("Z4").Value = "09:43"
$ws.Range("AB4").Value = "09:43"

# Row 5
$ws.Range("I5").NumberFormat = "@"
$ws.Range("A5").Value = 111803769
$ws.Range("B5").Value = 90666
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I5").Value = "10"
$ws.Range("Z5").Value = "09:59"
$ws.Range("AB5").Value = "09:59"

# Row 6
$ws.Range("I6").NumberFormat = "@"
$ws.Range("A6").Value = 111804210
$ws.Range("B6").Value = 85062
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 249278
$ws.Range("F6").Value = "Barrviolspindling"
$ws.Range("G6").Value = "Cortinarius harcynicus"
$ws.Range("H6").Value = "(Pers.) M.M.Moser"
$ws.Range("I6").Value = "1"
$ws.Range("Q6").Value = 695942.6774062206
$ws.Range("R6").Value = 6553663.74395387
$ws.Range("S6").Value = 1
$ws.Range("Z6").Value = "10:26"
$ws.Range("AB6").Value = "10:26"

# Row 7
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("A7").Value = 111891039
$ws.Range("B7").Value = 108219
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 219711
$ws.Range("F7").Value = "Sårläka"
$ws.Range("G7").Value = "Sanicula europaea"
$ws.Range("H7").Value = "L."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = "Ornö 1, Srm"
$ws.Range("Q7").Value = 695862.1592837617
$ws.Range("R7").Value = 6553543.752691799
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Stockholm"
$ws.Range("U7").Value = "Haninge"
$ws.Range("V7").Value = "Södermanland"
$ws.Range("W7").Value = "Ornö"
$ws.Range("Y7").Value = "2023-09-03"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").Value = "2023-09-03"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = ""
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Måns Persson"
$ws.Range("AX7").Value = "Måns Persson"
$ws.Range("AY7").Value = ""

